# Update column C ("Förändrad") date value from 46061 (2026-02-08) to
# 46062 (2026-02-09) for every data row on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRowByColumn = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162
$lastRowByUsed = $ws.UsedRange.Rows.Count
$lastRow = [Math]::Max($lastRowByColumn, $lastRowByUsed)

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
